$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -ne $null) {
        $text = $val.ToString()
        if ($text.Contains(",")) {
            $rawParts = $text.Split(",")
            $parts = @()
            foreach ($p in $rawParts) {
                $parts += $p.Trim()
            }

            $idx = [Array]::IndexOf($parts, "System")

            if ($idx -ge 0) {
                if ($idx -gt 0) {
                    $before = $parts[0..($idx-1)]
                } else {
                    $before = @()
                }
                if ($idx -lt ($parts.Length - 1)) {
                    $after = $parts[($idx+1)..($parts.Length-1)]
                } else {
                    $after = @()
                }
                $newParts = @($parts[$idx]) + $before + $after
            } else {
                $newParts = @()
                for ($i = $parts.Length - 1; $i -ge 0; $i--) {
                    $newParts += $parts[$i]
                }
            }

            $cell.Value2 = [string]::Join(", ", $newParts)
        }
    }
}
